$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 14, shifting existing rows 14-21 down to 15-22
# (mirrors Excel's Insert > Entire Row, with cells shifting down)
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new weekly price-report record
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C14").Value = "Arica y Parinacota"
$ws.Range("D14").Value = 44784
$ws.Range("E14").Value = 15
$ws.Range("F14").Value = 100112026
$ws.Range("G14").Value = "Haba"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 1000
$ws.Range("K14").Value = 1200
$ws.Range("L14").Value = 1300
$ws.Range("M14").Value = 1250
$ws.Range("N14").Value = "$/kilo"
$ws.Range("O14").Value = "Región de Arica y Parinacota"
$ws.Range("P14").Value = 1250
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = "Hortaliza"
